$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    (0.001, 20.17678),
    (0.001984596919383877, 64.87384525750433),
    (0.003868173634726945, 76.85628397922264),
    (0.006051410282056411, 81.26546465228964),
    (0.01001120224044809, 85.06828772318971),
    (0.01968593718743749, 91.26443002797724),
    (0.0250370074014803, 93.08310049749753),
    (0.02728445689137828, 95.1610403112026),
    (0.02959611922384477, 95.38924286600356),
    (0.03447629525905181, 97.38972944328131),
    (0.03650970194038808, 97.55952273279601),
    (0.04027685537107421, 98.9448572254126),
    (0.04209621924384876, 98.60316176060762),
    (0.04408681736347269, 99.71001229152353),
    (0.04449349869973995, 99.73581373282001),
    (0.04549949989997999, 99.65861097392623),
    (0.04599179835967193, 99.6895525034941),
    (0.04993018603720744, 100.7838227180577),
    (0.05042248449689938, 100.8236767534422),
    (0.05149269853970794, 100.7670245533888),
    (0.05196359271854371, 100.7899515108494),
    (0.05487457491498299, 101.6072723156725),
    (0.05513142628525705, 101.6300382324279),
    (0.05573074614922984, 101.6062818875768),
    (0.0559875975195039, 101.6216106563779),
    (0.05703640728145629, 101.9036968261809),
    (0.05776415283056611, 101.8831463939445),
    (0.05797819563912782, 101.8962756656992),
    (0.0589627925585117, 102.172521258124),
    (0.05988317663532707, 102.1385558048806),
    (0.06159551910382076, 103.2758473949729),
    (0.06324364872974594, 102.9763586476114),
    (0.06429245849169833, 103.0108921517697),
    (0.06502020404080816, 103.0038393554248),
    (0.06583356671334266, 103.0076617985294),
    (0.06660412082416484, 102.9810334044963),
    (0.0669251850370074, 102.9949831539166),
    (0.07092778555711142, 104.0970692338103),
    (0.07139867973594718, 104.1320418188029),
    (0.07227625525105021, 104.1011311515488),
    (0.07362472494498899, 104.117939449797),
    (0.07454510902180436, 104.077236326668),
    (0.0748875775155031, 104.0912044814027),
    (0.08019583916783356, 105.2317132624075),
    (0.08137307461492298, 105.2061595148771),
    (0.08235767153430686, 105.2143336021544),
    (0.08334226845369073, 105.2121522882839),
    (0.08770874174834967, 105.2131183810803),
    (0.08869333866773355, 105.2106788889844),
    (0.08967793558711742, 105.2199046177534),
    (0.09070534106821364, 105.1860151445006),
    (0.09248189637927585, 105.5426264095507),
    (0.09427985597119423, 105.1738623574632),
    (0.09575675135027005, 105.2915149021277),
    (0.09862492498499699, 103.7276165515901),
    (0.1007225445089018, 100.5547864362607),
    (0.108, 76.77531)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}